$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates to Price (D) and Volume(1h) (E) columns.
# Values are forced to text via a leading apostrophe in .Formula so that
# numeric-looking strings (e.g. '245.12', '1.0000') are NOT reinterpreted
# by Excel as numbers -- matching the original inlineStr text cells.
$ws.Range("D2").Formula = "'29.771.69"
$ws.Range("E2").Formula = "'  +2.00%  "
$ws.Range("D3").Formula = "'1.856.39"
$ws.Range("E3").Formula = "'  +1.60%  "
$ws.Range("E4").Formula = "'  +0.04%  "
$ws.Range("D5").Formula = "'245.12"
$ws.Range("E5").Formula = "'  +1.31%  "
$ws.Range("D6").Formula = "'0.6417"
$ws.Range("E6").Formula = "'  +3.33%  "
$ws.Range("D8").Formula = "'47.11"
$ws.Range("E8").Formula = "'  +4.76%  "
$ws.Range("D9").Formula = "'0.07500"
$ws.Range("E9").Formula = "'  +2.10%  "
$ws.Range("D10").Formula = "'0.2973"
$ws.Range("E10").Formula = "'  +2.33%  "
$ws.Range("D11").Formula = "'24.18"
$ws.Range("E11").Formula = "'  +4.83%  "
$ws.Range("D12").Formula = "'0.07680"
$ws.Range("E12").Formula = "'  -0.02%  "
$ws.Range("D13").Formula = "'1.872.12"
$ws.Range("E13").Formula = "'  +2.52%  "
$ws.Range("E14").Formula = "'  +1.94%  "
$ws.Range("D15").Formula = "'0.6853"
$ws.Range("E15").Formula = "'  +2.99%  "
$ws.Range("D16").Formula = "'83.99"
$ws.Range("E16").Formula = "'  +1.98%  "
$ws.Range("D17").Formula = "'0.000009510"
$ws.Range("E17").Formula = "'  +6.08%  "
$ws.Range("D18").Formula = "'6.088"
$ws.Range("E18").Formula = "'  +3.82%  "
$ws.Range("D19").Formula = "'29.751.14"
$ws.Range("E19").Formula = "'  +2.02%  "
$ws.Range("D20").Formula = "'2.122.63"
$ws.Range("E20").Formula = "'  +2.47%  "
$ws.Range("D21").Formula = "'239.98"
$ws.Range("E21").Formula = "'  +0.85%  "
$ws.Range("E22").Formula = "'  +1.81%  "
$ws.Range("D23").Formula = "'1.0000"
$ws.Range("E23").Formula = "'  -0.03%  "
$ws.Range("D24").Formula = "'7.443"
$ws.Range("E24").Formula = "'  +2.10%  "
$ws.Range("E25").Formula = "'  +0.06%  "
$ws.Range("E26").Formula = "'  +0.35%  "
$ws.Range("E27").Formula = "'  +0.36%  "
$ws.Range("D28").Formula = "'8.518"
$ws.Range("E28").Formula = "'  +0.35%  "
$ws.Range("D30").Formula = "'0.06178"
$ws.Range("E30").Formula = "'  +9.59%  "
$ws.Range("D31").Formula = "'1.499"
$ws.Range("E31").Formula = "'  +1.35%  "
$ws.Range("D32").Formula = "'1.274"
$ws.Range("E32").Formula = "'  +5.69%  "
$ws.Range("D33").Formula = "'4.145"
$ws.Range("E33").Formula = "'  +1.06%  "
$ws.Range("D34").Formula = "'4.111"
$ws.Range("E34").Formula = "'  +0.46%  "
$ws.Range("D35").Formula = "'1.879"
$ws.Range("E35").Formula = "'  +1.94%  "
$ws.Range("D36").Formula = "'1.159"
$ws.Range("E36").Formula = "'  +2.38%  "
$ws.Range("D37").Formula = "'0.7359"
$ws.Range("E37").Formula = "'  -0.11%  "
$ws.Range("E38").Formula = "'  -0.76%  "
$ws.Range("D39").Formula = "'2.857"
$ws.Range("E39").Formula = "'  +0.43%  "
$ws.Range("D40").Formula = "'0.01795"
$ws.Range("E40").Formula = "'  +1.49%  "
$ws.Range("D41").Formula = "'1.212.98"
$ws.Range("E41").Formula = "'  -0.17%  "
$ws.Range("D42").Formula = "'0.9227"
$ws.Range("E42").Formula = "'  +0.25%  "
$ws.Range("D43").Formula = "'6.186"
$ws.Range("E43").Formula = "'  -2.03%  "
$ws.Range("E44").Formula = "'  +0.00%  "
$ws.Range("D45").Formula = "'2.025.98"
$ws.Range("E45").Formula = "'  +2.62%  "
$ws.Range("D46").Formula = "'102.15"
$ws.Range("E46").Formula = "'  +0.32%  "
$ws.Range("D47").Formula = "'66.34"
$ws.Range("E47").Formula = "'  +2.10%  "
$ws.Range("D48").Formula = "'0.00000000122"
$ws.Range("E48").Formula = "'  +3.83%  "
$ws.Range("D49").Formula = "'9.331"
$ws.Range("E49").Formula = "'  +2.21%  "
$ws.Range("D50").Formula = "'0.4083"
$ws.Range("E50").Formula = "'  +1.14%  "

# Row 51: coin replaced from Algorand to Cronos (name, link, price, volume).
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Formula = "'0.05803"
$ws.Range("E51").Formula = "'  +0.69%  "
